# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.630.50"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "3.163.44"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.09"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.97"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.543"
$ws.Range("E8").Value = "  +15.94%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +6.69%  "
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").Value = "3.708.47"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.78"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D16").Value = "58.664.56"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").Value = "3.178.38"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.96"
$ws.Range("E20").Value = "  +4.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.12"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +5.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.73"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +12.98%  "
$ws.Range("D28").Value = "0.0₃0870"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.36"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.89"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.06"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.17"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("E34").Value = "  +4.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.95"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  +5.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.11"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0696"
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "2.666.23"
$ws.Range("E40").Value = "  +7.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.28"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.10"
$ws.Range("E43").Value = "  +4.21%  "
$ws.Range("E44").Value = "  +7.85%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "3.203.90"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  +14.87%  "
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.13"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.752"
$ws.Range("E51").Value = "  +1.95%  "
